# Apply scheduled market-price / profit refresh values across all Leve sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1454.4642
$ws.Range("I15").Value = 1454.4642
$ws.Range("K15").Value = 4363.392599999999
$ws.Range("M15").Value = -4194.392599999999
$ws.Range("H17").Value = 4281.1
$ws.Range("J17").Value = 3368.5
$ws.Range("L17").Value = 10105.5
$ws.Range("N17").Value = -10441.5
$ws.Range("H43").Value = 996.5833
$ws.Range("I43").Value = 996.5833
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 996.5833
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -927.5833
$ws.Range("N43").ClearContents()
$ws.Range("H62").Value = 4531.1665
$ws.Range("I62").Value = 4567.3335
$ws.Range("K62").Value = 4567.3335
$ws.Range("M62").Value = -3943.3335
$ws.Range("H65").Value = 4531.1665
$ws.Range("I65").Value = 4567.3335
$ws.Range("K65").Value = 22836.6675
$ws.Range("M65").Value = -19716.6675
$ws.Range("H96").Value = 421.05884
$ws.Range("J96").Value = 232.22223
$ws.Range("L96").Value = 696.66669
$ws.Range("N96").Value = -3442.66669
$ws.Range("H101").Value = 255.6923
$ws.Range("I101").Value = 262.16666
$ws.Range("J101").Value = 250.14285
$ws.Range("K101").Value = 786.4999799999999
$ws.Range("L101").Value = 750.4285500000001
$ws.Range("M101").Value = 835.5000200000001
$ws.Range("N101").Value = -3994.42855
$ws.Range("H132").Value = 15492.979
$ws.Range("I132").Value = 4312.1025
$ws.Range("K132").Value = 12936.3075
$ws.Range("M132").Value = -10406.3075
$ws.Range("H137").Value = 4311.5776
$ws.Range("I137").Value = 4761.5586
$ws.Range("J137").Value = 2920.7273
$ws.Range("K137").Value = 14284.6758
$ws.Range("L137").Value = 8762.1819
$ws.Range("M137").Value = -11734.6758
$ws.Range("N137").Value = -13862.1819

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4183.5293
$ws.Range("I45").Value = 3685.5715
$ws.Range("J45").Value = 4532.1
$ws.Range("K45").Value = 3685.5715
$ws.Range("L45").Value = 4532.1
$ws.Range("M45").Value = -3308.5715
$ws.Range("N45").Value = -5286.1
$ws.Range("H61").Value = 1850.6
$ws.Range("I61").Value = 1225.2106
$ws.Range("K61").Value = 1225.2106
$ws.Range("M61").Value = -1013.2106
$ws.Range("H74").Value = 1784.3077
$ws.Range("I74").Value = 930.6
$ws.Range("K74").Value = 930.6
$ws.Range("M74").Value = -56.60000000000002
$ws.Range("H77").Value = 1784.3077
$ws.Range("I77").Value = 930.6
$ws.Range("K77").Value = 4653
$ws.Range("M77").Value = -285
$ws.Range("H132").Value = 1388.0454
$ws.Range("I132").Value = 1387.762
$ws.Range("K132").Value = 4163.286
$ws.Range("M132").Value = -1633.286
$ws.Range("H136").Value = 1850.6
$ws.Range("I136").Value = 1225.2106
$ws.Range("K136").Value = 3675.6318
$ws.Range("M136").Value = -1125.6318

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3840.4614
$ws.Range("I105").Value = 2716.8462
$ws.Range("J105").Value = 4964.077
$ws.Range("K105").Value = 2716.8462
$ws.Range("L105").Value = 4964.077
$ws.Range("M105").Value = -969.8462
$ws.Range("N105").Value = -8458.077000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4840.1875
$ws.Range("I16").Value = 5371.6924
$ws.Range("J16").Value = 2537
$ws.Range("K16").Value = 5371.6924
$ws.Range("L16").Value = 2537
$ws.Range("M16").Value = -5084.6924
$ws.Range("N16").Value = -3111
$ws.Range("H99").Value = 18704380
$ws.Range("I99").Value = 6099536
$ws.Range("J99").Value = 25006804
$ws.Range("K99").Value = 6099536
$ws.Range("L99").Value = 25006804
$ws.Range("M99").Value = -6098038
$ws.Range("N99").Value = -25009800
$ws.Range("H105").Value = 2966.647
$ws.Range("I105").Value = 1786.5834
$ws.Range("J105").Value = 5798.8
$ws.Range("K105").Value = 1786.5834
$ws.Range("L105").Value = 5798.8
$ws.Range("M105").Value = -39.58339999999998
$ws.Range("N105").Value = -9292.799999999999
$ws.Range("H113").Value = 4840.1875
$ws.Range("I113").Value = 5371.6924
$ws.Range("J113").Value = 2537
$ws.Range("K113").Value = 5371.6924
$ws.Range("L113").Value = 2537
$ws.Range("M113").Value = -3201.6924
$ws.Range("N113").Value = -6877
$ws.Range("H126").Value = 18704380
$ws.Range("I126").Value = 6099536
$ws.Range("J126").Value = 25006804
$ws.Range("K126").Value = 18298608
$ws.Range("L126").Value = 75020412
$ws.Range("M126").Value = -18296138
$ws.Range("N126").Value = -75025352
$ws.Range("H134").Value = 3119.9048
$ws.Range("I134").Value = 2904.3333
$ws.Range("K134").Value = 8712.999899999999
$ws.Range("M134").Value = -6177.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1777
$ws.Range("I5").Value = 1917.2307
$ws.Range("J5").Value = 1412.4
$ws.Range("K5").Value = 5751.6921
$ws.Range("L5").Value = 4237.200000000001
$ws.Range("M5").Value = -5639.6921
$ws.Range("N5").Value = -4461.200000000001
$ws.Range("H135").Value = 1777
$ws.Range("I135").Value = 1917.2307
$ws.Range("J135").Value = 1412.4
$ws.Range("K135").Value = 17255.0763
$ws.Range("L135").Value = 12711.6
$ws.Range("M135").Value = -14720.0763
$ws.Range("N135").Value = -17781.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5000
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5504
$ws.Range("H40").Value = 11795.4
$ws.Range("J40").Value = 11795.4
$ws.Range("L40").Value = 11795.4
$ws.Range("N40").Value = -12097.4
$ws.Range("H102").Value = 2382.96
$ws.Range("I102").Value = 1644.25
$ws.Range("J102").Value = 3696.2222
$ws.Range("K102").Value = 1644.25
$ws.Range("L102").Value = 3696.2222
$ws.Range("M102").Value = -22.25
$ws.Range("N102").Value = -6940.2222
$ws.Range("H123").Value = 41600
$ws.Range("J123").Value = 41600
$ws.Range("L123").Value = 41600
$ws.Range("N123").Value = -46500
$ws.Range("H132").Value = 2890.6177
$ws.Range("I132").Value = 2882.1035
$ws.Range("K132").Value = 8646.3105
$ws.Range("M132").Value = -6116.3105

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3353.6086
$ws.Range("I82").Value = 1992
$ws.Range("J82").Value = 4601.75
$ws.Range("K82").Value = 1992
$ws.Range("L82").Value = 4601.75
$ws.Range("M82").Value = -1631
$ws.Range("N82").Value = -5323.75
$ws.Range("H85").Value = 3353.6086
$ws.Range("I85").Value = 1992
$ws.Range("J85").Value = 4601.75
$ws.Range("K85").Value = 1992
$ws.Range("L85").Value = 4601.75
$ws.Range("M85").Value = -744
$ws.Range("N85").Value = -7097.75
$ws.Range("H122").Value = 6443.875
$ws.Range("I122").Value = 3279.8
$ws.Range("K122").Value = 9839.400000000001
$ws.Range("M122").Value = -7389.400000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12351123
$ws.Range("J81").Value = 22227922
$ws.Range("L81").Value = 44455844
$ws.Range("N81").Value = -44457966
$ws.Range("H84").Value = 12351123
$ws.Range("J84").Value = 22227922
$ws.Range("L84").Value = 222279220
$ws.Range("N84").Value = -222289828
$ws.Range("H126").Value = 2031.25
$ws.Range("J126").Value = 3862.25
$ws.Range("L126").Value = 11586.75
$ws.Range("N126").Value = -16526.75
$ws.Range("H132").Value = 2310.8064
$ws.Range("I132").Value = 1643.48
$ws.Range("K132").Value = 4930.440000000001
$ws.Range("M132").Value = -2400.440000000001
$ws.Range("H136").Value = 1337.1154
$ws.Range("I136").Value = 683.3
$ws.Range("K136").Value = 2049.9
$ws.Range("M136").Value = 500.1000000000004
